$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 22, pushing the old row 22 ("Nason Creek Lower 15")
# down to row 23. This makes room for a "new" row 22 ("Nason Creek Lower 14") that
# carries the data that used to live in row 21 before the renumbering below.
$ws.Rows.Item(22).Insert()

# ---------------------------------------------------------------------------
# Row 19: "Nason Creek Lower 12" -> "Nason Creek Lower 11" (EDT/Okanogan update)
# ---------------------------------------------------------------------------
$ws.Range("A19").Value = "Nason Creek Lower 11"
$ws.Range("J19").Value = 5
$ws.Range("K19").Value = 3
$ws.Range("M19").Value = 3
$ws.Range("O19").Value = 1
$ws.Range("Q19").Value = 3
$ws.Range("R19").Value = 2
$ws.Range("T19").Value = 22
$ws.Range("U19").Value = 0.4888888888888889
$ws.Range("X19").Value = "Off-Channel-Side-Channels,PoolQuantity&Quality,Temperature-Rearing"
$ws.Range("Y19").Value = "Stability,Cover-Wood,Flow-SummerBaseFlow,Off-Channel-Floodplain,Riparian"
$ws.Range("Z19").Value = "Stability,Cover-Wood,Flow-SummerBaseFlow,Off-Channel-Floodplain,Off-Channel-Side-Channels,PoolQuantity&Quality,Riparian,Temperature-Rearing"

# ---------------------------------------------------------------------------
# Row 20: "Nason Creek Lower 13" -> "Nason Creek Lower 12" (label only)
# ---------------------------------------------------------------------------
$ws.Range("A20").Value = "Nason Creek Lower 12"

# ---------------------------------------------------------------------------
# Row 21: "Nason Creek Lower 14" -> "Nason Creek Lower 13"
# ---------------------------------------------------------------------------
$ws.Range("A21").Value = "Nason Creek Lower 13"
$ws.Range("O21").Value = 3
$ws.Range("T21").Value = 17
$ws.Range("U21").Value = 0.3777777777777778
$ws.Range("X21").Value = "Cover-Wood,Off-Channel-Floodplain,Off-Channel-Side-Channels,Riparian,Temperature-Rearing"
$ws.Range("Y21").Value = "Stability,CoarseSubstrate,Flow-SummerBaseFlow,PoolQuantity&Quality"

# ---------------------------------------------------------------------------
# Row 22 (new row): "Nason Creek Lower 14" - carries the data that row 21
# ("Nason Creek Lower 14") used to hold before today's update.
# ---------------------------------------------------------------------------
$ws.Range("A22").Value = "Nason Creek Lower 14"
$ws.Range("B22").Value = "Wenatchee"
$ws.Range("C22").Value = "Lower Nason Creek"
$ws.Range("D22").Value = "yes"
$ws.Range("E22").Value = "yes"
$ws.Range("F22").Value = "yes"
$ws.Range("G22").Value = 3
$ws.Range("H22").Value = 3
$ws.Range("I22").Value = 3
$ws.Range("J22").Value = 3
$ws.Range("K22").Value = 1
$ws.Range("L22").Value = 3
$ws.Range("M22").Value = 1
$ws.Range("N22").Value = 1
$ws.Range("O22").Value = 1
$ws.Range("P22").Value = 1
$ws.Range("Q22").Value = 1
$ws.Range("R22").Value = 1
$ws.Range("S22").Value = 1
$ws.Range("T22").Value = 15
$ws.Range("U22").Value = 0.3333333333333333
$ws.Range("V22").Value = 5
$ws.Range("W22").Value = 1
$ws.Range("X22").Value = "Cover-Wood,Off-Channel-Floodplain,Off-Channel-Side-Channels,PoolQuantity&Quality,Riparian,Temperature-Rearing"
$ws.Range("Y22").Value = "Stability,CoarseSubstrate,Flow-SummerBaseFlow"
$ws.Range("Z22").Value = "Stability,CoarseSubstrate,Cover-Wood,Flow-SummerBaseFlow,Off-Channel-Floodplain,Off-Channel-Side-Channels,PoolQuantity&Quality,Riparian,Temperature-Rearing"

# Row 23 ("Nason Creek Lower 15") already holds the correct (unchanged) data
# after the row insert shifted it down from row 22 - nothing else to update.
